$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for columns C,D,E,F,G,H,J,K,M,O across rows 2-25
# (columns A, B, I, L, N are left unchanged).
$cols = @("C","D","E","F","G","H","J","K","M","O")

$data = @{
    2 = @(4.917506253870776, 9.696610877278069, 13.72316626336374, 29.6099356654798, 28.99315455355914, 14.44729006807796, 9.911994425744528, 14.10343815003689, 16.80414170420432, 22.00389371730394)
    3 = @(4.749967865506476, 9.674964650762169, 13.75201404084057, 29.74426335292127, 29.17473534243246, 14.51812966689989, 9.94597073408687, 13.50087353984994, 16.55297541014077, 22.13074518779434)
    4 = @(4.645491827508812, 9.662934505219589, 13.77213439681326, 29.8354222229928, 29.2983656463721, 14.56451092798892, 9.968094908805655, 13.11646922174864, 16.39811253292635, 22.21461342758883)
    5 = @(4.602589742871721, 9.658352420355696, 13.78093870230616, 29.8747445970914, 29.35177645697625, 14.58413707858025, 9.977428685119964, 12.95638012159332, 16.33490770073204, 22.25029046147701)
    6 = @(4.5954484803122, 9.657611005239222, 13.7824371855153, 29.88140512208815, 29.36082766573555, 14.58743979476639, 9.978997772776015, 12.92959557066955, 16.32440865780165, 22.25630509173292)
    7 = @(4.644914454135577, 9.662871408578141, 13.77225068524274, 29.83594374569102, 29.29907372231997, 14.5647726761651, 9.968219499148038, 13.11432387196642, 16.39726043601362, 22.21508851206016)
    8 = @(4.860116706718534, 9.688887624216692, 13.7326130748869, 29.65444500336345, 29.05323103635409, 14.47111658292362, 9.923447713486961, 13.89877238499454, 16.71771257518236, 22.04638840378065)
    9 = @(5.266315428261748, 9.749756310899205, 13.67400095835449, 29.36778117698427, 28.66841932506133, 14.31035938735783, 9.845645254196659, 15.31564081568551, 17.33817253108766, 21.76319992149926)
    10 = @(5.551448256862575, 9.800265929788482, 13.64260756302609, 29.1998871089558, 28.4463039870212, 14.20622291105895, 9.794544589754381, 16.27480739136481, 17.78550089420347, 21.58442199459504)
    11 = @(5.67762347969589, 9.824450889115779, 13.63086117659455, 29.13288550869847, 28.35868685349344, 14.16188527567567, 9.772606876965188, 16.69213334872774, 17.98640727841314, 21.50950282022686)
    12 = @(5.724848352590953, 9.833778192643861, 13.62677755532459, 29.10886976604905, 28.32746027365183, 14.14553253480087, 9.764487277979931, 16.84734719943774, 18.06205424757143, 21.48205895839591)
    13 = @(5.714703001592534, 9.831761952275651, 13.62764082685627, 29.11398153483243, 28.33409831590913, 14.1490349464169, 9.766227635292502, 16.81404568983853, 18.04578245922945, 21.48792820323989)
    14 = @(5.681520069614502, 9.82521489292543, 13.63051791021285, 29.13088249019315, 28.35607857828812, 14.16053116499996, 9.771935112319762, 16.70495981348923, 17.99263979174341, 21.50722640200742)
    15 = @(5.661120996713789, 9.821226493846355, 13.63232767138387, 29.14141168595827, 28.36979698211949, 14.16762985405621, 9.77545554438737, 16.63777221248384, 17.96003038372461, 21.51916789960461)
    16 = @(5.543127180475933, 9.798709293409393, 13.64342622276592, 29.20445517185181, 28.45230197580045, 14.20918159284022, 9.796004562693083, 16.24714410783157, 17.7723139654783, 21.58944750356406)
    17 = @(5.46980176278721, 9.785201799015644, 13.65088403723645, 29.24553728428597, 28.5063700569975, 14.23544987402913, 9.808945512469931, 16.00257361244745, 17.65645024195686, 21.6342064719067)
    18 = @(5.427296808108196, 9.777546631041069, 13.65541215182326, 29.27004849997472, 28.53873087909535, 14.25084428486542, 9.816511964255955, 15.86011857815327, 17.58956857372369, 21.66055347797962)
    19 = @(5.412850182257271, 9.774974440118271, 13.65698626771099, 29.27849878070093, 28.54990382451209, 14.25610560087159, 9.819094999563998, 15.81158209513892, 17.5668842564631, 21.66957750802776)
    20 = @(5.4776419012147, 9.786627931279089, 13.65006544799459, 29.24107269177662, 28.50048361543791, 14.2326240089517, 9.80755518252791, 16.02879391881899, 17.66880940729311, 21.62937937671487)
    21 = @(5.691282108919874, 9.827133375223298, 13.62966294965407, 29.12588139545295, 28.34956929229483, 14.15714258307005, 9.770253595868498, 16.73707808626785, 18.00826127009602, 21.50153287099635)
    22 = @(5.827655688693706, 9.854588680163404, 13.61845313969179, 29.05850645892332, 28.26232614928037, 14.11035834333358, 9.746968954378291, 17.18352430953831, 18.22756906762319, 21.42337992652231)
    23 = @(5.755182264618941, 9.839846957395904, 13.62424166753287, 29.09373937384439, 28.30784043887867, 14.13509469489359, 9.759296420188301, 16.94677835876874, 18.11077231026275, 21.46459567022701)
    24 = @(5.474098458701302, 9.785982832435554, 13.65043478300452, 29.24308835438814, 28.50314090150932, 14.23390067168867, 9.808183356536613, 16.01694546840755, 17.66322266525295, 21.63155979273275)
    25 = @(5.158522975462247, 9.732256547274149, 13.68780861768319, 29.43786477903777, 28.76197038491345, 14.35139573561624, 9.865626241652796, 14.94622151470824, 17.1715588155867, 21.83468575480094)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
